$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts the existing rows 42-50
# down to 43-51 (preserving their data and formatting).
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new weekly price record.
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Macroferia Regional de Talca"
$ws.Range("C42").Value = "Maule"
$ws.Range("D42").Value = 44505
$ws.Range("E42").Value = 7
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100107
$ws.Range("H42").Value = "Otros"
$ws.Range("I42").Value = 100107002
$ws.Range("J42").Value = "Chirimoya"
$ws.Range("K42").Value = "Cultivar IV Región"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 25000
$ws.Range("O42").Value = 25000
$ws.Range("P42").Value = 25000
$ws.Range("Q42").Value = "`$/bandeja 10 kilos"
$ws.Range("R42").Value = "Provincia de Limarí"
$ws.Range("S42").Value = 2500
$ws.Range("T42").Value = 10
